# Finished Week 13 logging
# Update the "Road" (R) row totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet ("OFF") - row 3 is the "R" (Road) row
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 219
$wsOff.Range("C3").Value = 158

# DEF sheet ("DEF") - row 3 is the "R" (Road) row
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 409
$wsDef.Range("C3").Value = 290
$wsDef.Range("D3").Value = 114
$wsDef.Range("E3").Value = 42
